$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 54604
$ws.Range("F4").Value = 1330
$ws.Range("F5").Value = 367
$ws.Range("F6").Value = 317
$ws.Range("F7").Value = 873
$ws.Range("F8").Value = 741
$ws.Range("F9").Value = 387
$ws.Range("F10").Value = 3031
$ws.Range("F11").Value = 896
$ws.Range("F12").Value = 5205
$ws.Range("F13").Value = 1274
$ws.Range("G13").Value = 85
$ws.Range("F14").Value = 985
$ws.Range("F16").Value = 771
$ws.Range("F18").Value = 394
$ws.Range("F19").Value = 1256
$ws.Range("F21").Value = 37
$ws.Range("F22").Value = 168
$ws.Range("F23").Value = 353
$ws.Range("F24").Value = 16
$ws.Range("F25").Value = 36
$ws.Range("F27").Value = 66
$ws.Range("F28").Value = 57
$ws.Range("F29").Value = 4917
$ws.Range("F30").Value = 36
$ws.Range("F31").Value = 4890
$ws.Range("F32").Value = 8860
$ws.Range("F33").Value = 112
$ws.Range("F35").Value = 131
$ws.Range("F36").Value = 216
$ws.Range("F37").Value = 421
$ws.Range("F38").Value = 109
$ws.Range("F39").Value = 81
$ws.Range("F40").Value = 4189
$ws.Range("F41").Value = 225

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4177
$ws.Range("F12").Value = 1124

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 783
$ws.Range("F3").Value = 567
$ws.Range("F5").Value = 36

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 783
$ws.Range("F3").Value = 567
$ws.Range("F5").Value = 1330
$ws.Range("F6").Value = 367
$ws.Range("F7").Value = 317
$ws.Range("F8").Value = 873
$ws.Range("F9").Value = 741
$ws.Range("F10").Value = 387
$ws.Range("F11").Value = 3031
$ws.Range("F12").Value = 896
$ws.Range("F14").Value = 1274
$ws.Range("G14").Value = 85
$ws.Range("F15").Value = 36
$ws.Range("F17").Value = 985
$ws.Range("F19").Value = 771
$ws.Range("F20").Value = 394
$ws.Range("F22").Value = 1256
$ws.Range("F25").Value = 168
$ws.Range("F27").Value = 353
$ws.Range("F28").Value = 36
$ws.Range("F29").Value = 66
$ws.Range("F30").Value = 57
$ws.Range("F31").Value = 4917
$ws.Range("F32").Value = 36
$ws.Range("F33").Value = 4891
$ws.Range("F34").Value = 8860
$ws.Range("F35").Value = 112
$ws.Range("F37").Value = 131
$ws.Range("F38").Value = 216
$ws.Range("F39").Value = 421
$ws.Range("F42").Value = 109
$ws.Range("F43").Value = 81
$ws.Range("F44").Value = 4189
$ws.Range("F47").Value = 225
